$d = $word.ActiveDocument

$replacements = @(
    @{old="48×64="; new="55×59="},
    @{old="26×17="; new="52×47="},
    @{old="84×55="; new="33×53="},
    @{old="47×92="; new="75×90="},
    @{old="16×79="; new="37×50="},
    @{old="98×98="; new="94×49="},
    @{old="48×30="; new="22×98="},
    @{old="31×22="; new="28×68="},
    @{old="86×41="; new="40×38="},
    @{old="23×86="; new="96×21="},
    @{old="14×11="; new="19×22="},
    @{old="72×40="; new="94×34="},
    @{old="38×66="; new="52×48="},
    @{old="70×17="; new="82×95="},
    @{old="81×44="; new="68×28="},
    @{old="94×84="; new="13×25="},
    @{old="76×47="; new="72×41="},
    @{old="94×21="; new="15×31="},
    @{old="11×55="; new="90×45="},
    @{old="75×30="; new="94×86="},
    @{old="23×77="; new="32×98="},
    @{old="42×67="; new="55×26="},
    @{old="90×33="; new="53×73="},
    @{old="34×91="; new="36×44="},
    @{old="96×28="; new="25×60="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
